# Updates cryptos list values (Price and Volume(1h) columns) per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Prefix with an apostrophe so Excel always stores the value as literal text
    # (prevents numeric/date auto-conversion of values such as "59.20" or "1.00"),
    # then reset the style so no stray "text format" style is left on the cell.
    $ws.Range($cell).Value = "`'" + $text
    $ws.Range($cell).Style = "Normal"
}

Set-TextValue "D2" "36.308.98"
Set-TextValue "E2" "  -1.30%  "
Set-TextValue "D3" "2.040.30"
Set-TextValue "E3" "  -1.84%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "244.47"
Set-TextValue "E5" "  -0.17%  "
Set-TextValue "E6" "  +1.13%  "
Set-TextValue "D8" "54.57"
Set-TextValue "E8" "  +0.39%  "
Set-TextValue "D9" "59.20"
Set-TextValue "E9" "  -0.76%  "
Set-TextValue "E10" "  -0.41%  "
Set-TextValue "D11" "0.0740"
Set-TextValue "E11" "  -2.70%  "
Set-TextValue "E12" "  -3.91%  "
Set-TextValue "D13" "0.911"
Set-TextValue "E13" "  +3.43%  "
Set-TextValue "D14" "14.32"
Set-TextValue "E14" "  -4.19%  "
Set-TextValue "D15" "2.340.91"
Set-TextValue "E15" "  -1.68%  "
Set-TextValue "D16" "5.34"
Set-TextValue "E16" "  -2.64%  "
Set-TextValue "D17" "2.041.36"
Set-TextValue "E17" "  -1.63%  "
Set-TextValue "D18" "17.47"
Set-TextValue "E18" "  +1.28%  "
Set-TextValue "D19" "36.103.14"
Set-TextValue "E19" "  -1.57%  "
Set-TextValue "D20" "71.26"
Set-TextValue "E20" "  -1.97%  "
Set-TextValue "D21" "0.0₃0853"
Set-TextValue "E21" "  -2.91%  "
Set-TextValue "D22" "236.02"
Set-TextValue "E22" "  -0.64%  "
Set-TextValue "D23" "5.18"
Set-TextValue "E23" "  -4.54%  "
Set-TextValue "E24" "  -0.01%  "
Set-TextValue "D25" "2.35"
Set-TextValue "E25" "  -2.62%  "
Set-TextValue "D26" "2.25"
Set-TextValue "E26" "  +4.39%  "
Set-TextValue "D27" "9.32"
Set-TextValue "E27" "  -5.04%  "
Set-TextValue "D28" "163.78"
Set-TextValue "E28" "  -2.05%  "
Set-TextValue "D29" "19.90"
Set-TextValue "E29" "  -3.43%  "
Set-TextValue "E30" "  -1.63%  "
Set-TextValue "E31" "  -1.58%  "
Set-TextValue "D32" "4.96"
Set-TextValue "E32" "  -5.64%  "
Set-TextValue "D33" "0.0597"
Set-TextValue "E33" "  -1.61%  "
Set-TextValue "E34" "  -6.62%  "
Set-TextValue "D35" "0.0898"
Set-TextValue "E35" "  +7.70%  "
Set-TextValue "D36" "0.999"
Set-TextValue "E36" "  -0.18%  "
Set-TextValue "D37" "1.82"
Set-TextValue "E37" "  -1.06%  "
Set-TextValue "D38" "2.21"
Set-TextValue "E38" "  -6.43%  "
Set-TextValue "D39" "5.03"
Set-TextValue "E39" "  +3.69%  "
Set-TextValue "E40" "  -5.61%  "
Set-TextValue "E41" "  +1.87%  "
Set-TextValue "E42" "  -2.82%  "
Set-TextValue "D43" "1.10"
Set-TextValue "E43" "  -4.60%  "
Set-TextValue "D44" "0.0904"
Set-TextValue "E44" "  -4.90%  "
Set-TextValue "D45" "92.53"
Set-TextValue "E45" "  -3.93%  "
Set-TextValue "D46" "1.396.99"
Set-TextValue "E46" "  +3.65%  "
Set-TextValue "D47" "7.48"
Set-TextValue "E47" "  +3.97%  "
Set-TextValue "D48" "15.48"
Set-TextValue "E48" "  -3.56%  "
Set-TextValue "D49" "2.94"
Set-TextValue "E49" "  +1.77%  "
Set-TextValue "E50" "  -6.95%  "
Set-TextValue "D51" "46.04"
Set-TextValue "E51" "  +1.74%  "
